{"js": "// Replace the date line and every three-digit x one-digit multiplication\n// problem in the practice table with the new day's values.\nconst replacements = [\n  [\"2024-09-04 Wednesday\", \"2024-09-05 Thursday\"],\n  [\"639\u00d78=5112\", \"263\u00d76=1578\"],\n  [\"370\u00d74=1480\", \"341\u00d73=1023\"],\n  [\"748\u00d72=1496\", \"964\u00d73=2892\"],\n  [\"583\u00d72=1166\", \"570\u00d79=5130\"],\n  [\"562\u00d75=2810\", \"746\u00d76=4476\"],\n  [\"506\u00d76=3036\", \"237\u00d78=1896\"],\n  [\"545\u00d75=2725\", \"638\u00d74=2552\"],\n  [\"332\u00d72=664\", \"710\u00d77=4970\"],\n  [\"794\u00d72=1588\", \"145\u00d72=290\"],\n  [\"605\u00d75=3025\", \"996\u00d74=3984\"],\n  [\"843\u00d76=5058\", \"835\u00d77=5845\"],\n  [\"359\u00d74=1436\", \"216\u00d74=864\"],\n  [\"202\u00d76=1212\", \"891\u00d76=5346\"],\n  [\"866\u00d72=1732\", \"824\u00d76=4944\"],\n  [\"915\u00d77=6405\", \"864\u00d73=2592\"],\n  [\"201\u00d79=1809\", \"890\u00d76=5340\"],\n  [\"489\u00d77=3423\", \"899\u00d79=8091\"],\n  [\"909\u00d79=8181\", \"390\u00d72=780\"],\n  [\"606\u00d74=2424\", \"103\u00d78=824\"],\n  [\"102\u00d79=918\", \"281\u00d79=2529\"],\n  [\"655\u00d74=2620\", \"264\u00d77=1848\"],\n  [\"474\u00d75=2370\", \"398\u00d75=1990\"],\n  [\"554\u00d72=1108\", \"245\u00d77=1715\"],\n  [\"216\u00d72=432\", \"585\u00d79=5265\"],\n  [\"487\u00d76=2922\", \"978\u00d72=1956\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-09-04 Wednesday\"; New = \"2024-09-05 Thursday\" },\n    @{ Old = \"639\u00d78=5112\"; New = \"263\u00d76=1578\" },\n    @{ Old = \"370\u00d74=1480\"; New = \"341\u00d73=1023\" },\n    @{ Old = \"748\u00d72=1496\"; New = \"964\u00d73=2892\" },\n    @{ Old = \"583\u00d72=1166\"; New = \"570\u00d79=5130\" },\n    @{ Old = \"562\u00d75=2810\"; New = \"746\u00d76=4476\" },\n    @{ Old = \"506\u00d76=3036\"; New = \"237\u00d78=1896\" },\n    @{ Old = \"545\u00d75=2725\"; New = \"638\u00d74=2552\" },\n    @{ Old = \"332\u00d72=664\";  New = \"710\u00d77=4970\" },\n    @{ Old = \"794\u00d72=1588\"; New = \"145\u00d72=290\" },\n    @{ Old = \"605\u00d75=3025\"; New = \"996\u00d74=3984\" },\n    @{ Old = \"843\u00d76=5058\"; New = \"835\u00d77=5845\" },\n    @{ Old = \"359\u00d74=1436\"; New = \"216\u00d74=864\" },\n    @{ Old = \"202\u00d76=1212\"; New = \"891\u00d76=5346\" },\n    @{ Old = \"866\u00d72=1732\"; New = \"824\u00d76=4944\" },\n    @{ Old = \"915\u00d77=6405\"; New = \"864\u00d73=2592\" },\n    @{ Old = \"201\u00d79=1809\"; New = \"890\u00d76=5340\" },\n    @{ Old = \"489\u00d77=3423\"; New = \"899\u00d79=8091\" },\n    @{ Old = \"909\u00d79=8181\"; New = \"390\u00d72=780\" },\n    @{ Old = \"606\u00d74=2424\"; New = \"103\u00d78=824\" },\n    @{ Old = \"102\u00d79=918\";  New = \"281\u00d79=2529\" },\n    @{ Old = \"655\u00d74=2620\"; New = \"264\u00d77=1848\" },\n    @{ Old = \"474\u00d75=2370\"; New = \"398\u00d75=1990\" },\n    @{ Old = \"554\u00d72=1108\"; New = \"245\u00d77=1715\" },\n    @{ Old = \"216\u00d72=432\";  New = \"585\u00d79=5265\" },\n    @{ Old = \"487\u00d76=2922\"; New = \"978\u00d72=1956\" }\n)\n\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $rep.Old\n    $find.Replacement.Text = $rep.New\n    $found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Text not found: $($rep.Old)\"\n    }\n}\n"}
